$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "iPhone4" row (row 2) is moved further down the list, past the
# iPhone5/6/6plus block and the blank spacer row, to sit just above the
# "iPad Mini" row. Concretely:
#   - row 2 (iPhone4/960/640) content is removed in place (no shifting,
#     so iPhone5/iPhone6/iPhone6plus/blank keep their row numbers)
#   - a new row is inserted above the old "iPad Mini" row (row 7),
#     shifting iPad Mini / ipad Other / iPad Pro down by one row
#   - the new row 7 is filled with the iPhone4 data

# Remember the iPhone4 row's values before clearing them (.Value has a
# quirky getter in this host; .Value2 reads back the real scalar).
$name4 = $ws.Range("B2").Value2
$w4 = $ws.Range("C2").Value2
$h4 = $ws.Range("D2").Value2

# Clear out the old iPhone4 row in place (no shifting), removing the row
# entirely (content + formatting) rather than leaving an empty shell.
$ws.Rows.Item(2).Clear()

# Insert a fresh row just above "iPad Mini" (currently row 7), shifting
# iPad Mini and the rows below it down by one. The new row naturally
# inherits the number format/style from the row below it.
$ws.Rows.Item(7).Insert()

# Fill the newly inserted row 7 with the iPhone4 data (now positioned
# right before iPad Mini).
$ws.Range("B7").Value = $name4
$ws.Range("C7").Value = $w4
$ws.Range("D7").Value = $h4
$ws.Range("E7").Formula = "=D7/C7"
$ws.Range("F7").Formula = "=D7/C7"

# Update the active selection/view: select K12 (moves the view so the
# previous topLeftCell pin is no longer needed).
$ws.Range("K12").Select()
